$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-label the rscKpi id values with zero-padded numbers
# (GP1/GP2/GP3 -> GP01/GP02/GP03, BP1/BP2/BP3 -> BP01/BP02/BP03)
$ws.Range("B16").Value = "GP01"
$ws.Range("B17:B18").Value = "GP02"
$ws.Range("B19:B20").Value = "GP03"
$ws.Range("B21:B25").Value = "BP01"
$ws.Range("B26:B30").Value = "BP02"
$ws.Range("B31:B60").Value = "BP03"

# Reset the view: scroll/select back to A1 instead of the previously
# saved G1/I4 position (frozen header rows stay frozen).
$ws.Range("A1").Select()
